# Fruta / hortaliza, semanal
#
# The sheet contains one data row per market report (rows 2-21). This edit
# re-shuffles the report rows: the "identity" columns of each market
# (Mercado ID, Mercado, Region, Codreg, Categoria ID, Categoria,
# Clasificacion -> A,B,C,E,F,G,R) stay put, while the report-specific
# columns (Fecha, Variedad, Calidad, Volumen, Precio minimo/maximo/promedio,
# Unidad de comercializacion, Origen, Precio $/Kg, Kg o Unidades ->
# D,H,I,J,K,L,M,N,O,P,Q) are redistributed across the rows according to a
# fixed permutation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 21

# Columns whose contents move between rows.
$cols = @("D", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q")

# Snapshot the current ("before") values of the movable columns for every
# data row, so the writes below don't clobber values we still need to read.
$original = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $values = @()
    foreach ($col in $cols) {
        $values += $ws.Range("$col$r").Value2
    }
    $original[$r] = $values
}

# Target row -> source row (the target row receives the original values
# that used to live in the source row).
$mapping = @{
    2  = 19
    3  = 10
    4  = 11
    5  = 6
    6  = 7
    7  = 8
    8  = 14
    9  = 18
    10 = 5
    11 = 12
    12 = 13
    13 = 17
    14 = 20
    15 = 21
    16 = 15
    17 = 16
    18 = 9
    19 = 2
    20 = 3
    21 = 4
}

for ($targetRow = $firstDataRow; $targetRow -le $lastDataRow; $targetRow++) {
    $sourceRow = $mapping[$targetRow]
    $values = $original[$sourceRow]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$targetRow").Value = $values[$i]
    }
}
